# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interest count) figures in column F on the
# 展览 (Exhibitions) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (rows keyed by their row number on this sheet) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 20
$ws1.Range("F3").Value = 8121
$ws1.Range("F5").Value = 944
$ws1.Range("F6").Value = 301
$ws1.Range("F7").Value = 809
$ws1.Range("F8").Value = 626
$ws1.Range("F9").Value = 105
$ws1.Range("F10").Value = 71
$ws1.Range("F13").Value = 3342
$ws1.Range("F15").Value = 114
$ws1.Range("F16").Value = 758
$ws1.Range("F17").Value = 761
$ws1.Range("F19").Value = 466
$ws1.Range("F21").Value = 288
$ws1.Range("F22").Value = 813
$ws1.Range("F23").Value = 364
$ws1.Range("F26").Value = 131
$ws1.Range("F27").Value = 297
$ws1.Range("F29").Value = 82
$ws1.Range("F31").Value = 507
$ws1.Range("F32").Value = 590
$ws1.Range("F33").Value = 29
$ws1.Range("F34").Value = 39
$ws1.Range("F35").Value = 22
$ws1.Range("F38").Value = 113

# --- Sheet: 全部类型 (same events, offset by the extra rows on this sheet) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 20
$ws4.Range("F5").Value = 8121
$ws4.Range("F7").Value = 944
$ws4.Range("F8").Value = 301
$ws4.Range("F9").Value = 809
$ws4.Range("F10").Value = 626
$ws4.Range("F11").Value = 105
$ws4.Range("F12").Value = 71
$ws4.Range("F16").Value = 3342
$ws4.Range("F18").Value = 114
$ws4.Range("F20").Value = 758
$ws4.Range("F21").Value = 761
$ws4.Range("F24").Value = 466
$ws4.Range("F26").Value = 288
$ws4.Range("F27").Value = 815
$ws4.Range("F28").Value = 364
$ws4.Range("F31").Value = 131
$ws4.Range("F32").Value = 297
$ws4.Range("F34").Value = 85
$ws4.Range("F36").Value = 507
$ws4.Range("F37").Value = 590
$ws4.Range("F38").Value = 29
$ws4.Range("F39").Value = 39
$ws4.Range("F40").Value = 22
$ws4.Range("F43").Value = 113
